# AddCandidatePO, excel operations and Test
# Re-create the "Candidate" sheet so it gets a fresh internal sheetId (matches
# the author's workbook having recreated this tab), while keeping its
# original position between "JD" and "Skills".
$wb = $excel.ActiveWorkbook

$newCandidate = $wb.Worksheets.Add()
$oldCandidate = $wb.Worksheets.Item("Candidate")
$oldCandidate.Delete()
$newCandidate.Name = "Candidate"
$skills = $wb.Worksheets.Item("Skills")
$newCandidate.Move($skills)

$ws = $wb.Worksheets.Item("Candidate")

# Header row - written in the same order the source workbook's shared
# strings table was built in (D,G,H,I filled before E,F).
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "LastName"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "contact_no"
$ws.Range("G1").Value = "current_company"
$ws.Range("H1").Value = "current_CTC"
$ws.Range("I1").Value = "expected_CTC"
$ws.Range("E1").Value = "Notice"
$ws.Range("F1").Value = "currency"

# Data row
$ws.Range("A2").Value = "Test_can_aut01_FN"
$ws.Range("B2").Value = "Test_can_aut01_LN"
$ws.Range("C2").Value = "test@0001.gmail"
$ws.Range("D2").Value = 943001022
$ws.Range("E2").Value = "Immediate"
$ws.Range("F2").Value = "USD*"
$ws.Range("G2").Value = "Infinite"
$ws.Range("H2").Value = 200000
$ws.Range("I2").Value = 300000

# Email hyperlink on the candidate's email cell
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:test@0001.gmail")

# Column widths to fit the new content
$ws.Columns.Item(3).ColumnWidth = 14.81640625
$ws.Columns.Item(4).ColumnWidth = 9.81640625

# The "JD" sheet keeps its own last selection, just moved from K2 to J2
$jd = $wb.Worksheets.Item("JD")
$jd.Activate()
$jd.Range("J2").Select()

# Make "Candidate" the active/selected sheet and cell (activated last so it
# ends up the workbook's active tab)
$ws.Activate()
$ws.Range("D2").Select()
